# "Add unix tools cert" commit - the visible OOXML change is the removal
# of the old "Note: please use the web link instead of this PDF" paragraph
# that used to sit right under the top horizontal rule.
$d = $word.ActiveDocument

$marker = "*Please consider using the following link for this CV*"

for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like $marker) {
        # Delete the whole paragraph, including its end-of-paragraph mark,
        # so no blank paragraph is left behind.
        $p.Range.Delete()
    }
}
